$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 837, shifting existing rows 837-878 down to 838-879.
$ws.Rows.Item(837).Insert()

# Populate the newly inserted row with the new day's data.
# Prefix the date with an apostrophe so Excel stores it as literal text
# (matching the existing "2026/..." cells) instead of auto-converting it
# to a date serial value; then reset the style so no quote-prefix style
# index lingers on the cell.
$ws.Range("A837").Value = "'2026/02/18"
$ws.Range("A837").Style = "Normal"
$ws.Range("B837").Value = "水"
$ws.Range("C837").Value = 7
$ws.Range("D837").Value = 201
